# "Updated policies and graphs"
#
# 1. The "Weights" row (row 7) changes the weight of the
#    "Temp Policy_Industries" column (AC) from 1 to 0, and the total
#    weight (AG7) drops from 13 to 12.
# 2. Because every LockdownEffectiveness value (column AG, rows 9-233)
#    is a weighted average over the Weights row, that reweighting
#    changes the cached AG value for every row whose weighted sum
#    actually includes column AC (i.e. every row from 24 down).
# 3. Twelve new daily rows (9/30/2020 - 10/11/2020) are appended at the
#    bottom of the "Converted Data" sheet, rows 222-233, following the
#    exact same pattern as the preceding rows (and the same new AG
#    weighted average, 0.25).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Converted Data")

# --- 1. Update the Weights row -------------------------------------------
$ws.Range("AC7").Value = 0
$ws.Range("AG7").Value = 12

# --- 2. Recomputed LockdownEffectiveness (column AG) for existing rows ---
$ws.Range("AG24:AG26").Value = 0.2750000000000001
$ws.Range("AG27:AG28").Value = 0.5250000000166667
$ws.Range("AG29:AG31").Value = 0.5527777777916668
$ws.Range("AG32:AG63").Value = 0.8333333333416665
$ws.Range("AG64:AG84").Value = 0.9166666666749999
$ws.Range("AG85:AG91").Value = 0.9000000000083332
$ws.Range("AG92:AG100").Value = 0.8444444444583333
$ws.Range("AG101:AG101").Value = 0.7611111111250001
$ws.Range("AG102:AG108").Value = 0.291666666675
$ws.Range("AG109:AG221").Value = 0.25

# --- 3. Append the 12 new daily rows (222-233) ----------------------------
$newDates = @("9/30/2020","10/1/2020","10/2/2020","10/3/2020","10/4/2020","10/5/2020","10/6/2020","10/7/2020","10/8/2020","10/9/2020","10/10/2020","10/11/2020")
$rowVals  = @(0,0,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,0,0)
$firstNewRow = 222
$lastNewRow  = $firstNewRow + $newDates.Length - 1

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $firstNewRow + $i
    # Leading apostrophe forces the date-looking text to stay literal text,
    # matching every other cell in column A (which are strings, not dates).
    $ws.Cells.Item($r, 1).Value = "'" + $newDates[$i]

    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, 2 + $c).Value = $rowVals[$c]
    }

    $ws.Cells.Item($r, 33).Value = 0.25
}

# Copy the formatting (bold font + border) of the last pre-existing date
# cell onto the new date cells so column A keeps a uniform look.
$ws.Range("A221").Copy()
$ws.Range("A" + $firstNewRow + ":A" + $lastNewRow).PasteSpecial(-4122)
